# Applies the cryptos-list refresh described in the commit:
# "Updated cryptos list on Mon Feb 12 21:32:05 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'50.100.30"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +4.41%  "
$ws.Range("D3").Value = "'2.638.30"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +5.67%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'327.67"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.31%  "
$ws.Range("D6").Value = "'110.96"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.35%  "
$ws.Range("E7").Value = "  +1.89%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  +4.66%  "
$ws.Range("D10").Value = "'40.76"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.39%  "
$ws.Range("D11").Value = "'20.71"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.50%  "
$ws.Range("D12").Value = "'0.0821"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.16%  "
$ws.Range("E13").Value = "  +0.84%  "
$ws.Range("D14").Value = "'7.30"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.93%  "
$ws.Range("D15").Value = "'3.059.66"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.97%  "
$ws.Range("D16").Value = "'2.642.17"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.86%  "
$ws.Range("D17").Value = "'0.881"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +5.39%  "
$ws.Range("D18").Value = "'50.052.95"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.60%  "
$ws.Range("D19").Value = "'3.07"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +11.51%  "
$ws.Range("D20").Value = "'13.30"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.12%  "
$ws.Range("D21").Value = "'6.83"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.75%  "
$ws.Range("D22").Value = "'0.0₃0963"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.70%  "
$ws.Range("D23").Value = "'73.05"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.21%  "
$ws.Range("D24").Value = "'279.87"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.74%  "
$ws.Range("D25").Value = "'2.60"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.73%  "
$ws.Range("D26").Value = "'26.65"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.19%  "
$ws.Range("D27").Value = "'0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.07%  "
$ws.Range("E28").Value = "  -0.86%  "
$ws.Range("B29").Value = "Kaspa"
$ws.Range("C29").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D29").Value = "'0.145"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.70%  "
$ws.Range("D30").Value = "'9.95"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.41%  "
$ws.Range("B31").Value = "InjectiveProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D31").Value = "'36.59"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.70%  "
$ws.Range("D32").Value = "'49.81"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.79%  "
$ws.Range("D33").Value = "'19.82"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.84%  "
$ws.Range("D34").Value = "'5.45"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.13%  "
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("D36").Value = "'0.0797"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.42%  "
$ws.Range("D37").Value = "'2.06"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +6.57%  "
$ws.Range("D38").Value = "'4.76"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.37%  "
$ws.Range("D39").Value = "'3.10"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +7.29%  "
$ws.Range("E40").Value = "  +1.54%  "
$ws.Range("D41").Value = "'123.50"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.36%  "
$ws.Range("D42").Value = "'22.51"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.35%  "
$ws.Range("E43").Value = "  +0.43%  "
$ws.Range("D44").Value = "'0.0314"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.95%  "
$ws.Range("D45").Value = "'3.34"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +6.36%  "
$ws.Range("D46").Value = "'2.065.60"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.95%  "
$ws.Range("D47").Value = "'2.30"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +15.45%  "
$ws.Range("D48").Value = "'2.00"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +8.97%  "
$ws.Range("D49").Value = "'9.07"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.19%  "
$ws.Range("E50").Value = "  +4.68%  "
$ws.Range("D51").Value = "'81.77"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.28%  "

Write-Output "Applied 95 cell updates"